# ---------------------------------------------------------------------------
# Task: Insert a new "Task 6" section (heading + 15 numbered steps) after the
# "10. Get the text..." paragraph (end of Task 5), move the _GoBack bookmark
# onto the new heading, and tidy the trailing blank paragraphs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# 1. Locate the paragraph that ends Task 5 ("...Please accept Terms and
#    Conditions..."). The new Task 6 content is inserted right after it.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
  if ($d.Paragraphs($i).Range.Text -like "*Please accept Terms and Conditions*") {
    $targetIndex = $i
  }
}

# 2. The paragraph immediately after it is the first of a run of blank
#    paragraphs at the end of the document. Mint a brand-new numbered-list
#    definition (decimal "1." list) on it first -- this allocates a fresh
#    w:numId (26 in this document) in numbering.xml that our new list
#    paragraphs below can reference directly.
$mintPara = $d.Paragraphs($targetIndex + 1)
$template = $word.ListGalleries(1).ListTemplates(1)
$mintPara.Range.ListFormat.ApplyListTemplateWithLevel($template)
$newNumId = $mintPara.Range.ListFormat.List.ListID

# 3. Insert the Task 6 heading + the 15 numbered steps as literal OOXML
#    right at the start of that blank paragraph. The first new paragraph
#    reuses/replaces that blank paragraph; the rest are added after it, so
#    the remaining trailing blank paragraphs are left undisturbed.
$insertPoint = $mintPara.Range
$insertPoint.Collapse(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:rStyle w:val="SubtleEmphasis"/><w:b/><w:bCs/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_Toc127286337"/><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rStyle w:val="SubtleEmphasis"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Task </w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rStyle w:val="SubtleEmphasis"/><w:b/><w:bCs/></w:rPr><w:t>6</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Navigate onto http://demo.openemr.io/b/openemr/</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Update username as admin</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Update password as pass</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Select language as English (Indian)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Click on the login button</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Click on Patient </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t> Click New Search</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Add the first name, last name</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Update DOB as today''s date driver.findElement(By.id("form_DOB")).sendKeys("2021-12-");</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Update the gender</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>. Click on the create new patient button above the form</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>. Click on confirm create new patient button.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>. Print the text from alert box (if any error before handling alert add 5 sec wait)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>. Handle alert</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Close the Happy Birthday popup</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="26"/></w:numPr><w:spacing w:after="0" w:line="254" w:lineRule="auto"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Get the added patient name and print in the console.</w:t></w:r></w:p>'
$insertPoint.InsertXML($xml)

# 4. The engine drops the run-level w:rStyle when InsertXML parses raw
#    OOXML, so re-apply the "Subtle Emphasis" character style to the two
#    runs of the new heading paragraph ("Task " and "6") explicitly.
$newCount = $d.Paragraphs.Count
$headingIndex = -1
for ($i = 1; $i -le $newCount; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text -like "*Task 6*" -and $p.Style.NameLocal -like "*Heading 3*") {
    $headingIndex = $i
  }
}
$headingPara = $d.Paragraphs($headingIndex)
$hr = $headingPara.Range
$headingTextRange = $d.Range($hr.Start, $hr.End - 1)
$headingTextRange.Style = "Subtle Emphasis"

# 5. Remove the old _GoBack bookmark that used to sit on the last paragraph
#    of the document -- it has now moved to the new Task 6 heading -- by
#    clearing that paragraph back to a clean empty paragraph.
$finalCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($finalCount)
$lastPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>')

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count
